# Remove the first 13 data rows (rows 2-14) from the sheet.
# Excel shifts the remaining rows (old 15..161) up to become new rows 2..148,
# and the used-range dimension shrinks from A1:R161 to A1:R148 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A14").EntireRow.Delete()
